$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: "News page fintune" task moves from "Not Started" to "In Progress"
# and gets a start date recorded.
$ws.Range("D16").Value = "In Progress"
$ws.Range("E16").Value = "21-07-2025"

# Row 17: new task "Author page" (High priority, Completed) started/ended 23-07-2025
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Author page "
$ws.Range("C17").Value = "High"
$ws.Range("D17").Value = "Completed"
$ws.Range("E17").Value = "23-07-2025"
$ws.Range("F17").Value = "23-07-2025"

# Row 18: new task "Sitemap update" (Low priority, Not Started) started/ended 23-07-2025
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Sitemap update"
$ws.Range("C18").Value = "Low"
$ws.Range("D18").Value = "Not Started"
$ws.Range("E18").Value = "23-07-2025"
$ws.Range("F18").Value = "23-07-2025"

# Row 23's height settles back to the standard 19.5 after editing nearby rows
$ws.Rows.Item(23).RowHeight = 19.5
